$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style (borders, font, number format, alignment) from A38 to A39 first,
# so the new date cell keeps the same formatting as the rest of column A.
$ws.Range("A38").Copy()
$ws.Range("A39").PasteSpecial(-4122)  # xlPasteFormats

# Append a new row of forecast data (row 39)
$ws.Range("A39").Value = 45986
$ws.Range("B39").Value = 2025
$ws.Range("C39").Value = 0.2194017515915414
$ws.Range("D39").Value = 2026
$ws.Range("E39").Value = -0.1883185981439661
